$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the shared-string value used by H1 from "BiomassToLAI" to "LeafBiomassToLAI"
$ws.Range("H1").Value = "LeafBiomassToLAI"

# Update the view: H1 becomes the active/selected cell
$ws.Range("H1").Select()
